$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.754521
$ws.Range("H2").Value = 2.263563
$ws.Range("I2").Value = 0.2768403531129761
$ws.Range("J2").Value = 0.2768403531129761
$ws.Range("M2").Value = 3.704480666666667
$ws.Range("N2").Value = 11.113442
$ws.Range("O2").Value = 0.298964201395561
$ws.Range("P2").Value = 0.2989642013955609
$ws.Range("Q2").Value = 2.795108457094
$ws.Range("R2").Value = 25.155976113846
$ws.Range("S2").Value = 0.082765355082486
$ws.Range("T2").Value = 0.08276535508248599
$ws.Range("G3").Value = 0.754521
$ws.Range("H3").Value = 2.263563
$ws.Range("I3").Value = 0.2768403531129761
$ws.Range("J3").Value = 0.2768403531129761
$ws.Range("O3").Value = 0.3765624616238499
$ws.Range("P3").Value = 0.3765624616238499
$ws.Range("Q3").Value = 3.520598507098
$ws.Range("R3").Value = 31.685386563882
$ws.Range("S3").Value = 0.1042476848450381
$ws.Range("T3").Value = 0.1042476848450381
$ws.Range("G4").Value = 0.754521
$ws.Range("H4").Value = 2.263563
$ws.Range("I4").Value = 0.2768403531129761
$ws.Range("J4").Value = 0.2768403531129761
$ws.Range("M4").Value = 2.870093333333334
$ws.Range("N4").Value = 8.610280000000001
$ws.Range("O4").Value = 0.2316263029934534
$ws.Range("P4").Value = 0.2316263029934534
$ws.Range("Q4").Value = 2.16554569196
$ws.Range("R4").Value = 19.48991122764
$ws.Range("S4").Value = 0.06412350751096084
$ws.Range("T4").Value = 0.06412350751096084
$ws.Range("G5").Value = 0.754521
$ws.Range("H5").Value = 2.263563
$ws.Range("I5").Value = 0.2768403531129761
$ws.Range("J5").Value = 0.2768403531129761
$ws.Range("M5").Value = 1.150472333333333
$ws.Range("N5").Value = 3.451417
$ws.Range("O5").Value = 0.09284703398713583
$ws.Range("P5").Value = 0.09284703398713583
$ws.Range("Q5").Value = 0.8680555354189999
$ws.Range("R5").Value = 7.812499818771
$ws.Range("S5").Value = 0.02570380567449118
$ws.Range("T5").Value = 0.02570380567449118
$ws.Range("G6").Value = 0.9731926666666667
$ws.Range("I6").Value = 0.3570728998754956
$ws.Range("J6").Value = 0.3570728998754956
$ws.Range("M6").Value = 3.704480666666667
$ws.Range("N6").Value = 11.113442
$ws.Range("O6").Value = 0.298964201395561
$ws.Range("P6").Value = 0.2989642013955609
$ws.Range("Q6").Value = 3.605173418608445
$ws.Range("R6").Value = 32.446560767476
$ws.Range("S6").Value = 0.1067520143512746
$ws.Range("T6").Value = 0.1067520143512746
$ws.Range("G7").Value = 0.9731926666666667
$ws.Range("I7").Value = 0.3570728998754956
$ws.Range("J7").Value = 0.3570728998754956
$ws.Range("O7").Value = 0.3765624616238499
$ws.Range("P7").Value = 0.3765624616238499
$ws.Range("R7").Value = 40.868293718092
$ws.Range("S7").Value = 0.1344602501562831
$ws.Range("T7").Value = 0.1344602501562831
$ws.Range("G8").Value = 0.9731926666666667
$ws.Range("I8").Value = 0.3570728998754956
$ws.Range("J8").Value = 0.3570728998754956
$ws.Range("M8").Value = 2.870093333333334
$ws.Range("N8").Value = 8.610280000000001
$ws.Range("O8").Value = 0.2316263029934534
$ws.Range("P8").Value = 0.2316263029934534
$ws.Range("Q8").Value = 2.793153784648889
$ws.Range("R8").Value = 25.13838406184
$ws.Range("S8").Value = 0.0827074756973126
$ws.Range("T8").Value = 0.0827074756973126
$ws.Range("G9").Value = 0.9731926666666667
$ws.Range("I9").Value = 0.3570728998754956
$ws.Range("J9").Value = 0.3570728998754956
$ws.Range("M9").Value = 1.150472333333333
$ws.Range("N9").Value = 3.451417
$ws.Range("O9").Value = 0.09284703398713583
$ws.Range("P9").Value = 0.09284703398713583
$ws.Range("Q9").Value = 1.119631238002889
$ws.Range("R9").Value = 10.076681142026
$ws.Range("S9").Value = 0.03315315967062529
$ws.Range("T9").Value = 0.03315315967062529
$ws.Range("G10").Value = 0.7824410000000001
$ws.Range("H10").Value = 2.347323
$ws.Range("I10").Value = 0.2870844452706686
$ws.Range("J10").Value = 0.2870844452706686
$ws.Range("M10").Value = 3.704480666666667
$ws.Range("N10").Value = 11.113442
$ws.Range("O10").Value = 0.298964201395561
$ws.Range("P10").Value = 0.2989642013955609
$ws.Range("Q10").Value = 2.898537557307333
$ws.Range("R10").Value = 26.086838015766
$ws.Range("S10").Value = 0.08582797191343305
$ws.Range("T10").Value = 0.08582797191343305
$ws.Range("G11").Value = 0.7824410000000001
$ws.Range("H11").Value = 2.347323
$ws.Range("I11").Value = 0.2870844452706686
$ws.Range("J11").Value = 0.2870844452706686
$ws.Range("O11").Value = 0.3765624616238499
$ws.Range("P11").Value = 0.3765624616238499
$ws.Range("Q11").Value = 3.650873357391334
$ws.Range("R11").Value = 32.85786021652201
$ws.Range("S11").Value = 0.1081052254050404
$ws.Range("T11").Value = 0.1081052254050404
$ws.Range("G12").Value = 0.7824410000000001
$ws.Range("H12").Value = 2.347323
$ws.Range("I12").Value = 0.2870844452706686
$ws.Range("J12").Value = 0.2870844452706686
$ws.Range("M12").Value = 2.870093333333334
$ws.Range("N12").Value = 8.610280000000001
$ws.Range("O12").Value = 0.2316263029934534
$ws.Range("P12").Value = 0.2316263029934534
$ws.Range("Q12").Value = 2.245678697826667
$ws.Range("R12").Value = 20.21110828044
$ws.Range("S12").Value = 0.06649630870497136
$ws.Range("T12").Value = 0.06649630870497138
$ws.Range("G13").Value = 0.7824410000000001
$ws.Range("H13").Value = 2.347323
$ws.Range("I13").Value = 0.2870844452706686
$ws.Range("J13").Value = 0.2870844452706686
$ws.Range("M13").Value = 1.150472333333333
$ws.Range("N13").Value = 3.451417
$ws.Range("O13").Value = 0.09284703398713583
$ws.Range("P13").Value = 0.09284703398713583
$ws.Range("Q13").Value = 0.9001767229656668
$ws.Range("R13").Value = 8.101590506691002
$ws.Range("S13").Value = 0.0266549392472238
$ws.Range("T13").Value = 0.0266549392472238
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2153186666666667
$ws.Range("H14").Value = 0.645956
$ws.Range("I14").Value = 0.07900230174085969
$ws.Range("J14").Value = 0.07900230174085969
$ws.Range("M14").Value = 3.704480666666667
$ws.Range("N14").Value = 11.113442
$ws.Range("O14").Value = 0.298964201395561
$ws.Range("P14").Value = 0.2989642013955609
$ws.Range("Q14").Value = 0.7976438378391111
$ws.Range("R14").Value = 7.178794540552
$ws.Range("S14").Value = 0.02361886004836725
$ws.Range("T14").Value = 0.02361886004836725
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2153186666666667
$ws.Range("H15").Value = 0.645956
$ws.Range("I15").Value = 0.07900230174085969
$ws.Range("J15").Value = 0.07900230174085969
$ws.Range("O15").Value = 0.3765624616238499
$ws.Range("P15").Value = 0.3765624616238499
$ws.Range("Q15").Value = 1.004677903487111
$ws.Range("R15").Value = 9.042101131384001
$ws.Range("S15").Value = 0.02974930121748829
$ws.Range("T15").Value = 0.02974930121748829
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2153186666666667
$ws.Range("H16").Value = 0.645956
$ws.Range("I16").Value = 0.07900230174085969
$ws.Range("J16").Value = 0.07900230174085969
$ws.Range("M16").Value = 2.870093333333334
$ws.Range("N16").Value = 8.610280000000001
$ws.Range("O16").Value = 0.2316263029934534
$ws.Range("P16").Value = 0.2316263029934534
$ws.Range("Q16").Value = 0.6179846697422222
$ws.Range("R16").Value = 5.561862027680001
$ws.Range("S16").Value = 0.0182990110802086
$ws.Range("T16").Value = 0.0182990110802086
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.2153186666666667
$ws.Range("H17").Value = 0.645956
$ws.Range("I17").Value = 0.07900230174085969
$ws.Range("J17").Value = 0.07900230174085969
$ws.Range("M17").Value = 1.150472333333333
$ws.Range("N17").Value = 3.451417
$ws.Range("O17").Value = 0.09284703398713583
$ws.Range("P17").Value = 0.09284703398713583
$ws.Range("Q17").Value = 0.2477181688502222
$ws.Range("R17").Value = 2.229463519652
$ws.Range("S17").Value = 0.007335129394795559
$ws.Range("T17").Value = 0.007335129394795559
